$d = $word.ActiveDocument
$HIGHLIGHT_COLOR = 5258796  # BGR for hex 2C3E50 (w:color w:val="2C3E50")

# Applies bold + color to the first occurrence of $searchText found at or after
# $startPos within the given paragraph's range. Returns the end position of the
# matched (and now-formatted) run so subsequent calls can continue the scan
# forward from there.
function Highlight-Next($paragraph, $startPos, $searchText) {
    $r = $paragraph.Range.Duplicate
    $r.Start = $startPos
    $ok = $r.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Text not found: $searchText"
    }
    $r.Font.Bold = 1
    $r.Font.Color = $HIGHLIGHT_COLOR
    return $r.End
}

# Bolds/colors every term in $terms, in order, within $paragraph.
function Highlight-Terms($paragraph, $terms) {
    $pos = $paragraph.Range.Start
    foreach ($term in $terms) {
        $pos = Highlight-Next $paragraph $pos $term
    }
}

# --- Professional Experience bullets (Siege Analytics) ---------------------

Highlight-Terms ($d.Paragraphs.Item(10)) @("23%", "64%")

Highlight-Terms ($d.Paragraphs.Item(12)) @("±4.2%", "±2.1%", "71%", "87%")

Highlight-Terms ($d.Paragraphs.Item(13)) @("73.5%", "$4.7M")

Highlight-Terms ($d.Paragraphs.Item(14)) @("$2")

# --- Key Achievements and Impact bullets ------------------------------------

Highlight-Terms ($d.Paragraphs.Item(50)) @("73.5%")

Highlight-Terms ($d.Paragraphs.Item(51)) @("$4.7M")

Highlight-Terms ($d.Paragraphs.Item(53)) @("178%")

Write-Output "done"
